# Fill in the "Start"/"End" (first shift) time columns for the first nine
# weekly blocks of the timesheet (weeks starting at rows 2, 11, 20, 29, 38,
# 47, 56, 65 and the partial week at row 74), and remove the stray sample
# times that had been left in row 146.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$weekStarts = @(2, 11, 20, 29, 38, 47, 56, 65, 74)

foreach ($start in $weekStarts) {
    if ($start -eq 74) {
        $end = 77
    } else {
        $end = $start + 6
    }

    $cRange = $ws.Range("C$start`:C$end")
    $cRange.Value = 0.0625
    $cRange.NumberFormat = "h:mm:ss"

    $dRange = $ws.Range("D$start`:D$end")
    $dRange.Value = 0.1979166666666667
    $dRange.NumberFormat = "h:mm:ss"
}

# Row 146 previously carried one-off sample Start/End values in C:F;
# clear them out entirely (removing the cells, not just their contents).
$ws.Range("C146:F146").Clear()
